$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — match formatting of existing header H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-29.
$data = @{
    2  = @(1, 5)
    3  = @(1, 7)
    4  = @(1, 5)
    5  = @(1, 6)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(4, 5)
    9  = @(1, 3)
    10 = @(1, 7)
    11 = @(1, 3)
    12 = @(1, 5)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 4)
    19 = @(1, 3)
    20 = @(1, 6)
    21 = @(1, 5)
    22 = @(1, 5)
    23 = @(1, 4)
    24 = @(1, 4)
    25 = @(4, 6)
    26 = @(1, 2)
    27 = @(1, 3)
    28 = @(1, 2)
    29 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
